$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$count = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldText) {
        $cell.Value2 = $newText
        $count++
    }
}
Write-Host "changed=$count"
